$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.29663348197937
$ws.Range("B1").Value = 1.931168675422668
$ws.Range("C1").Value = 5.261541843414307
$ws.Range("D1").Value = 1.931468367576599
$ws.Range("E1").Value = 1.092854022979736
